$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-12-28 Thursday" "2023-12-29 Friday"

Replace-Text "31×63=" "94×81="
Replace-Text "54×95=" "46×40="
Replace-Text "99×32=" "30×90="
Replace-Text "13×42=" "66×99="
Replace-Text "51×32=" "31×35="
Replace-Text "73×93=" "91×68="
Replace-Text "25×42=" "68×90="
Replace-Text "69×55=" "26×82="
Replace-Text "84×58=" "61×78="
Replace-Text "47×72=" "92×83="
Replace-Text "29×12=" "26×64="
Replace-Text "60×24=" "12×27="
Replace-Text "43×33=" "70×34="
Replace-Text "55×52=" "94×18="
Replace-Text "12×17=" "56×50="
Replace-Text "31×54=" "88×48="
Replace-Text "36×70=" "87×71="
Replace-Text "15×88=" "88×37="
Replace-Text "47×49=" "96×44="
Replace-Text "77×65=" "55×64="
Replace-Text "11×61=" "60×57="
Replace-Text "92×90=" "55×28="
Replace-Text "59×89=" "31×66="
Replace-Text "26×27=" "48×84="
Replace-Text "70×14=" "96×49="
